$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.983.08"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("D3").Value = "1.819.39"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.57%  "

$ws.Range("D5").Value = "311.44"
$ws.Range("E5").Value = "  -1.02%  "

$ws.Range("E6").Value = "  -0.55%  "

$ws.Range("D7").Value = "0.4534"
$ws.Range("E7").Value = "  +6.39%  "

$ws.Range("D8").Value = "0.3697"
$ws.Range("E8").Value = "  +0.73%  "

$ws.Range("D9").Value = "0.07272"
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("D10").Value = "0.8538"
$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("D11").Value = "20.73"
$ws.Range("E11").Value = "  -1.02%  "

$ws.Range("D12").Value = "1.815.94"
$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("D13").Value = "6.638"
$ws.Range("E13").Value = "  -0.47%  "

$ws.Range("D14").Value = "92.34"
$ws.Range("E14").Value = "  +4.59%  "

$ws.Range("D15").Value = "0.07098"
$ws.Range("E15").Value = "  -0.38%  "

$ws.Range("D16").Value = "5.316"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").Value = "0.000008783"
$ws.Range("E18").Value = "  -0.84%  "

$ws.Range("E19").Value = "  -0.54%  "

$ws.Range("D20").Value = "14.94"
$ws.Range("E20").Value = "  -0.82%  "

$ws.Range("D21").Value = "26.997.62"
$ws.Range("E21").Value = "  -0.95%  "

$ws.Range("D22").Value = "5.159"
$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("D23").Value = "10.92"
$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").Value = "1.987"
$ws.Range("E24").Value = "  -0.87%  "

$ws.Range("D25").Value = "151.59"
$ws.Range("E25").Value = "  -1.37%  "

$ws.Range("D26").Value = "2.205"
$ws.Range("E26").Value = "  +4.32%  "

$ws.Range("E27").Value = "  +0.73%  "

$ws.Range("D28").Value = "5.233"
$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("D29").Value = "116.26"
$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("D30").Value = "0.08856"
$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("D31").Value = "1.183"
$ws.Range("E31").Value = "  -1.00%  "

$ws.Range("D32").Value = "0.7501"
$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("D33").Value = "'2.940"
$ws.Range("E33").Value = "  +4.56%  "

$ws.Range("D34").Value = "4.434"
$ws.Range("E34").Value = "  -0.04%  "

$ws.Range("E35").Value = "  -0.60%  "

$ws.Range("D36").Value = "1.099"
$ws.Range("E36").Value = "  -1.32%  "

$ws.Range("D37").Value = "0.01962"
$ws.Range("E37").Value = "  -0.25%  "

$ws.Range("D38").Value = "0.05229"
$ws.Range("E38").Value = "  -0.81%  "

$ws.Range("D39").Value = "0.5297"
$ws.Range("E39").Value = "  +5.23%  "

$ws.Range("D40").Value = "7.176"
$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("D41").Value = "2.865"
$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("D42").Value = "0.1706"
$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("D43").Value = "0.5198"
$ws.Range("E43").Value = "  +9.92%  "

$ws.Range("D44").Value = "8.506"
$ws.Range("E44").Value = "  -1.26%  "

$ws.Range("D45").Value = "10.63"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("D46").Value = "1.953"
$ws.Range("E46").Value = "  +8.31%  "

$ws.Range("D47").Value = "105.35"
$ws.Range("E47").Value = "  -1.41%  "

$ws.Range("E48").Value = "  -0.60%  "

$ws.Range("D49").Value = "1.664"
$ws.Range("E49").Value = "  +0.46%  "

$ws.Range("D50").Value = "0.06373"
$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "0.9163"
$ws.Range("E51").Value = "  +0.11%  "
